$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Authors" field (column E) for this reference had previously been
# re-appended to itself on every re-run of the SLR importer, growing by one
# extra space between entries each time. Eventually the field became too
# large and broke the whole-sheet update. The fix re-runs the importer one
# more time, producing the next (24-space) generation of the field, and also
# adds a new "Misc. Data" column (J) that the importer now emits alongside
# the existing columns.

$authors = @(
    'Ke%Hu%NULL%1',
    'Wei-jie%Guan%NULL%1',
    'Ying%Bi%NULL%1',
    'Wei%Zhang%NULL%0',
    'Lanjuan%Li%NULL%0',
    'Boli%Zhang%NULL%1',
    'Qingquan%Liu%NULL%1',
    'Yuanlin%Song%NULL%1',
    'Xingwang%Li%NULL%0',
    'Zhongping%Duan%NULL%1',
    'Qingshan%Zheng%NULL%1',
    'Zifeng%Yang%NULL%1',
    'Jingyi%Liang%NULL%1',
    'Mingfeng%Han%NULL%0',
    'Lianguo%Ruan%NULL%1',
    'Chaomin%Wu%NULL%1',
    'Yunting%Zhang%NULL%1',
    'Zhen-hua%Jia%NULL%1',
    'Nan-shan%Zhong%NULL%0'
)

$separator = ',' + "".PadLeft(24, ' ')
$authorsField = '[' + ($authors -join $separator) + ']'

# New "Misc. Data" header column.
$ws.Range("J1").Value = "Misc. Data"

# Regenerated (longer) Authors field for the one oversized reference.
$ws.Range("E2").Value = $authorsField

# The reference has no misc. data recorded, so the new column is blank.
$ws.Range("J2").Value = ""
